$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Copy A1's current formatting (bold font / border / centered) onto B1:E1
#        while its old value ("Year") is still in place -- this captures the
#        exact same cell style (xf) used for the header row, with no new
#        style entries created. ---
$ws.Range("A1").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2) The new table only has 3 rows; drop the old rows 4 and 5 entirely
#        (this removes old A4/B4/A5 along with their row elements). ---
$ws.Rows("4:5").Delete()

# --- 3) The old A2/A3 cells (value + header-style) are not part of the new
#        table at all, so clear them completely (content + formatting). ---
$ws.Range("A2:A3").Clear()

# --- 4) Stage the new text values off to the side first. Typing a string
#        that "looks like" a number (e.g. "101,057", "2019") straight into a
#        cell makes Excel auto-convert it to a real number with a matching
#        number format, which would both change the stored type and mint a
#        new cell style. Driving it through a formula ("="101,057"") and then
#        converting that formula to its static result via copy/paste-values
#        keeps it as plain text (shared string) with the cell's original
#        "General" formatting untouched. ---
$stage = $ws.Range("Z1:Z9")
$stage.Cells.Item(1,1).Formula = '="Unnamed: 0"'
$stage.Cells.Item(2,1).Formula = '="2019"'
$stage.Cells.Item(3,1).Formula = '="Unnamed: 1"'
$stage.Cells.Item(4,1).Formula = '="2018"'
$stage.Cells.Item(5,1).Formula = '="Unnamed: 2"'
$stage.Cells.Item(6,1).Formula = '="101,057"'
$stage.Cells.Item(7,1).Formula = '="22"'
$stage.Cells.Item(8,1).Formula = '="96,293"'
$stage.Cells.Item(9,1).Formula = '="20"'

$stage.Copy()
$ws.Range("Z11").PasteSpecial(-4163)     # xlPasteValues
$excel.CutCopyMode = $false

# --- 5) Move the staged, plain-text values onto the real table cells. Using
#        paste-values here only changes the cell content, leaving whatever
#        formatting is already sitting on the destination cell (e.g. the
#        header style applied in step 1) exactly as it is. ---
$ws.Range("Z11").Copy(); $ws.Range("A1").PasteSpecial(-4163); $excel.CutCopyMode = $false
$ws.Range("Z12").Copy(); $ws.Range("B1").PasteSpecial(-4163); $excel.CutCopyMode = $false
$ws.Range("Z13").Copy(); $ws.Range("C1").PasteSpecial(-4163); $excel.CutCopyMode = $false
$ws.Range("Z14").Copy(); $ws.Range("D1").PasteSpecial(-4163); $excel.CutCopyMode = $false
$ws.Range("Z15").Copy(); $ws.Range("E1").PasteSpecial(-4163); $excel.CutCopyMode = $false
$ws.Range("Z16").Copy(); $ws.Range("B2").PasteSpecial(-4163); $excel.CutCopyMode = $false
$ws.Range("Z17").Copy(); $ws.Range("B3").PasteSpecial(-4163); $excel.CutCopyMode = $false
$ws.Range("Z18").Copy(); $ws.Range("D2").PasteSpecial(-4163); $excel.CutCopyMode = $false
$ws.Range("Z19").Copy(); $ws.Range("D3").PasteSpecial(-4163); $excel.CutCopyMode = $false

# --- 6) Remove the scratch/staging cells; nothing of them should remain. ---
$ws.Range("Z1:Z19").Clear()
